$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# hardware sheet: new rows / cells for additional hardware options
# ---------------------------------------------------------------------------
$hw = $wb.Worksheets.Item("hardware")

# Row 1 (headers) - add RPM column header in K1
$hw.Range("K1").Value = "RPM"

# Row 2 (Honeywell HR16 / RW) - add RPM value
$hw.Range("K2").Value = 6000

# Row 3 (Honeywell HR12 / RW) - add RPM value, drop the volume-total formula
$hw.Range("K3").Value = 6000
$hw.Range("L3").ClearContents()

# Row 4 (was "ATA " IMU placeholder) -> now "Honeywell MIMU" IMU with real values
$hw.Range("A4").Value = "Honeywell MIMU"
$hw.Range("E4").Value = 4.7
$hw.Range("F4").Value = 32
$hw.Range("G4").Value = 32
$hw.Range("H4").Formula = "=0.169*PI()*0.233^2"
$hw.Range("L4").ClearContents()

# Row 5 (A-STR / star) - drop the total-volume formula
$hw.Range("L5").ClearContents()

# Row 6 (NEW): Honeywell HR0610 reaction wheel
$hw.Range("A6").Value = "Honeywell HR0610"
$hw.Range("B6").Value = "RW"
$hw.Range("C6").Value = 12
$hw.Range("D6").Value = 0.055
$hw.Range("E6").Value = 15
$hw.Range("F6").Value = 80
$hw.Range("G6").Value = 15
$hw.Range("H6").Formula = "=0.265*0.267*0.12"
$hw.Range("I6").Value = 4
$hw.Range("K6").Value = 6000

# Row 7 (NEW): Honeywell HR14 reaction wheel
$hw.Range("A7").Value = "Honeywell HR14"
$hw.Range("B7").Value = "RW"
$hw.Range("C7").Value = 75
$hw.Range("D7").Value = 0.2

# Row 8 (NEW): MSL Terminal descent sensor (TDS)
$hw.Range("A8").Value = "MSL Terminal descent"
$hw.Range("B8").Value = "TDS"
$hw.Range("E8").Value = 25
$hw.Range("G8").Value = 30
$hw.Range("H8").Formula = "=1.3*0.5*0.4"
$hw.Range("L8").Value = "Deep Space Comm book p455"

$hw.Columns.Item(1).ColumnWidth = 20.57

# ---------------------------------------------------------------------------
# orb_mission sheet: additional approach parameters
# ---------------------------------------------------------------------------
$om = $wb.Worksheets.Item("orb_mission")

$om.Range("A10").Value = "app_slew"
$om.Range("B10").Value = 90
$om.Range("C10").Value = "deg"

$om.Range("A11").Value = "app_time"
$om.Range("B11").Value = 120
$om.Range("C11").Value = "sec"

# ---------------------------------------------------------------------------
# probe_props sheet: more thrusters + lander/heatshield geometry
# ---------------------------------------------------------------------------
$pp = $wb.Worksheets.Item("probe_props")

$pp.Range("B3").Value = 8

$pp.Range("A5").Value = "d_lander"
$pp.Range("B5").Value = 1.2
$pp.Range("C5").Value = "m"
$pp.Range("D5").Value = "ref sub_output"

$pp.Range("A6").Value = "d_hs"
$pp.Range("B6").Formula = "=B5*1.7"
$pp.Range("C6").Value = "m"

$pp.Range("A7").Value = "h_lander"
$pp.Range("B7").Value = "ref sub_output"
$pp.Range("C7").Value = 0.72

$pp.Columns.Item(1).ColumnWidth = 12.57

# ---------------------------------------------------------------------------
# Selections (match final cursor positions per sheet). Do this before moving
# sheets around so each selection "follows" the correct logical worksheet.
# ---------------------------------------------------------------------------
$om.Range("B11").Select()
$pp.Range("D8").Select()

# ---------------------------------------------------------------------------
# Reorder sheets: orb_props should now come before probe_props
# ---------------------------------------------------------------------------
$op = $wb.Worksheets.Item("orb_props")
$op.Move($pp)

$hw.Activate()
$hw.Range("B12").Select()
